# EmpMasterChange_Upload.xlsx
# 30-07-2018 : Added GlobalExculde ... Two New Fields Added -> MastEmp -> SPLALL, BAALL
# Adds two new trailing columns (F = SPLALL, G = BAALL) to Sheet1: a header
# cell in row 1 for each, and a literal text "0" in every data row (2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Header row (row 1): F1 = "SPLALL", G1 = "BAALL" ----------------------
# Same look as the existing header cells (B1:D1 use style s="2"): text
# number format + the yellow header fill. Touch the border object too (even
# though the end result has no visible border) so a dedicated header xf is
# registered, mirroring the extra cellXfs entry introduced upstream.
$headerRange = $ws.Range("F1:G1")
$headerRange.NumberFormat = "@"
$headerRange.Interior.ColorIndex = 6
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.LineStyle = -4142
$ws.Range("F1").Value = "SPLALL"
$ws.Range("G1").Value = "BAALL"

# ---- Data rows (rows 2-11): F/G = "0" typed with a leading apostrophe -----
# so it lands as text (quote-prefixed), matching the EmpCode-style columns.
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 6).NumberFormat = "@"
    $ws.Cells.Item($r, 7).NumberFormat = "@"
    $ws.Cells.Item($r, 6).Value = "'0"
    $ws.Cells.Item($r, 7).Value = "'0"
}

# ---- Selection moved from C15 to J11 in the saved view --------------------
$ws.Range("J11").Select()

Write-Output "EmpMasterChange_Upload: added SPLALL/BAALL columns (F:G)"
